$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row-2 "id" cell (A2): shared string "F_01" ---
$ws.Range("A2").Value = "F_01"

# --- Updated fairytale description (B2): literal "<br>" markers inserted
#     after each sentence ---
$ws.Range("B2").Value = 'In the fairytale "Alice in Daschland," a curious young girl named Alice discovers a magical realm powered by Dasch, a company known for its repository for humanities data. <br>This world is filled with technological wonders like talking robots, digital trees, and houses that change color. <br>Guided by a robot named RabbIT, Alice explores various marvels, including a café with size-altering tea and a library with books that speak. <br>She meets the Queen of Hearts, an engineer who creates devices connecting emotions with digital data. <br>Although Alice returns to her world with a token from Daschland, the experience leaves her inspired, reminding her that the realm of innovation and dreams, where humanities data come to life, awaits her return.'

# --- A2 gets its own (explicitly applied) font, distinct from the default style ---
$ws.Range("A2").Font.Name = "Arial"

# --- Column B sized to fit the long wrapped description; row 1 kept at the
#     sheet default height (explicit/custom), row 2 grown for the wrapped text ---
$ws.Columns("B").ColumnWidth = 93.66666666666667
$ws.Rows(1).RowHeight = 15.75
$ws.Rows(2).RowHeight = 98

# --- Wrap text on the description column header + value ---
$ws.Range("B1:B2").WrapText = $true

# --- View state: zoomed in, selection parked further down the sheet ---
$excel.ActiveWindow.Zoom = 150
$ws.Range("B27").Select()

Write-Host "done"
